$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H137").Value = 3093.5278
$ws.Range("I137").Value = 2584.9656
$ws.Range("J137").Value = 5200.4287
$ws.Range("K137").Value = 7754.8968
$ws.Range("L137").Value = 15601.2861
$ws.Range("M137").Value = -5204.8968
$ws.Range("N137").Value = -20701.2861

$ws.Range("H138").Value = 2712.9778
$ws.Range("I138").Value = 2414.5789
$ws.Range("J138").Value = 2931.0386
$ws.Range("K138").Value = 7243.736699999999
$ws.Range("L138").Value = 8793.1158
$ws.Range("M138").Value = -2103.736699999999
$ws.Range("N138").Value = -19073.1158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 501530.94
$ws.Range("I32").Value = 579664.7
$ws.Range("K32").Value = 579664.7
$ws.Range("M32").Value = -579377.7

$ws.Range("H88").Value = 2001.75
$ws.Range("J88").Value = 2001.75
$ws.Range("L88").Value = 2001.75
$ws.Range("N88").Value = -2813.75

$ws.Range("H91").Value = 2001.75
$ws.Range("J91").Value = 2001.75
$ws.Range("L91").Value = 2001.75
$ws.Range("N91").Value = -4809.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3420
$ws.Range("I86").Value = 3420
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3420
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2297
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3420
$ws.Range("I89").Value = 3420
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17100
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11484
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 3037.2083
$ws.Range("I134").Value = 2579.9375
$ws.Range("J134").Value = 3951.75
$ws.Range("K134").Value = 7739.8125
$ws.Range("L134").Value = 11855.25
$ws.Range("M134").Value = -5204.8125
$ws.Range("N134").Value = -16925.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 556.25
$ws.Range("I22").Value = 476.92307
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 476.92307
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -126.92307
$ws.Range("N22").Value = -1600

$ws.Range("H31").Value = 5091.35
$ws.Range("I31").Value = 1307.1515
$ws.Range("K31").Value = 1307.1515
$ws.Range("M31").Value = -1012.1515

$ws.Range("H34").Value = 5091.35
$ws.Range("I34").Value = 1307.1515
$ws.Range("K34").Value = 1307.1515
$ws.Range("M34").Value = -1105.1515

$ws.Range("H132").Value = 3969883.8
$ws.Range("I132").Value = 1273.3214
$ws.Range("J132").Value = 11907104
$ws.Range("K132").Value = 3819.9642
$ws.Range("L132").Value = 35721312
$ws.Range("M132").Value = -1289.9642
$ws.Range("N132").Value = -35726372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 105.71429
$ws.Range("K2").Value = 120
$ws.Range("L2").Value = 634.28574
$ws.Range("M2").Value = -7
$ws.Range("N2").Value = -860.28574

$ws.Range("H88").Value = 11016.1
$ws.Range("J88").Value = 11016.1
$ws.Range("L88").Value = 33048.3
$ws.Range("N88").Value = -33904.3

$ws.Range("H91").Value = 11016.1
$ws.Range("J91").Value = 11016.1
$ws.Range("L91").Value = 33048.3
$ws.Range("N91").Value = -36012.3

$ws.Range("H122").Value = 4711.64
$ws.Range("J122").Value = 11189.8
$ws.Range("L122").Value = 100708.2
$ws.Range("N122").Value = -105608.2

$ws.Range("H132").Value = 2351.7173
$ws.Range("I132").Value = 2475.4
$ws.Range("K132").Value = 22278.6
$ws.Range("M132").Value = -19748.6

$ws.Range("H136").Value = 6160
$ws.Range("I136").Value = 1800
$ws.Range("J136").Value = 9066.667
$ws.Range("K136").Value = 5400
$ws.Range("L136").Value = 27200.001
$ws.Range("M136").Value = -300
$ws.Range("N136").Value = -37400.001

$ws.Range("H137").Value = 6952719
$ws.Range("I137").Value = 27795066
$ws.Range("K137").Value = 83385198
$ws.Range("M137").Value = -83380098

$ws.Range("H139").Value = 4299.6787
$ws.Range("I139").Value = 1044.1666
$ws.Range("J139").Value = 6741.3125
$ws.Range("K139").Value = 3132.4998
$ws.Range("L139").Value = 20223.9375
$ws.Range("M139").Value = 2007.5002
$ws.Range("N139").Value = -30503.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5733.6113
$ws.Range("I70").Value = 5427
$ws.Range("J70").Value = 5978.9
$ws.Range("K70").Value = 5427
$ws.Range("L70").Value = 5978.9
$ws.Range("M70").Value = -5157
$ws.Range("N70").Value = -6518.9

$ws.Range("H73").Value = 5733.6113
$ws.Range("I73").Value = 5427
$ws.Range("J73").Value = 5978.9
$ws.Range("K73").Value = 5427
$ws.Range("L73").Value = 5978.9
$ws.Range("M73").Value = -4491
$ws.Range("N73").Value = -7850.9

$ws.Range("H97").Value = 2550.9092
$ws.Range("I97").Value = 2717.7778
$ws.Range("K97").Value = 2717.7778
$ws.Range("M97").Value = -2221.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 100000
$ws.Range("J81").Value = 100000
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996

$ws.Range("H82").Value = 62502724
$ws.Range("I82").Value = 83336460
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 83336460
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -83336099
$ws.Range("N82").Value = -2222

$ws.Range("H84").Value = 100000
$ws.Range("J84").Value = 100000
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984

$ws.Range("H85").Value = 62502724
$ws.Range("I85").Value = 83336460
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 83336460
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -83335212
$ws.Range("N85").Value = -3996

$ws.Range("H136").Value = 5209581
$ws.Range("I136").Value = 1071.1538
$ws.Range("J136").Value = 27779792
$ws.Range("K136").Value = 3213.4614
$ws.Range("L136").Value = 83339376
$ws.Range("M136").Value = -663.4614000000001
$ws.Range("N136").Value = -83344476

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 112247336
$ws.Range("I2").Value = 10000000
$ws.Range("J2").Value = 125028250
$ws.Range("K2").Value = 10000000
$ws.Range("L2").Value = 125028250
$ws.Range("M2").Value = -9999888
$ws.Range("N2").Value = -125028474

$ws.Range("H23").Value = 1200
$ws.Range("I23").Value = 1200
$ws.Range("K23").Value = 1200
$ws.Range("M23").Value = -971

$ws.Range("H62").Value = 8823.333
$ws.Range("J62").Value = 9988
$ws.Range("L62").Value = 9988
$ws.Range("N62").Value = -11236

$ws.Range("H65").Value = 8823.333
$ws.Range("J65").Value = 9988
$ws.Range("L65").Value = 49940
$ws.Range("N65").Value = -56180

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 4631941
$ws.Range("I132").Value = 2591.8572
$ws.Range("J132").Value = 11113030
$ws.Range("K132").Value = 7775.571599999999
$ws.Range("L132").Value = 33339090
$ws.Range("M132").Value = -5245.571599999999
$ws.Range("N132").Value = -33344150

$ws.Range("H136").Value = 1771.069
$ws.Range("I136").Value = 1339.1224
$ws.Range("J136").Value = 4122.778
$ws.Range("K136").Value = 4017.3672
$ws.Range("L136").Value = 12368.334
$ws.Range("M136").Value = -1467.3672
$ws.Range("N136").Value = -17468.334
